$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the f6bf9a95 row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-10-27 06:24:15"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the f6bf9a95 row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-10-27 06:24:03"
$wsZhCn.Range("K4").Value = "2016-10-27 06:24:44"

# de-de sheet: Correspond Handback DateTime for the f6bf9a95 row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-10-27 06:25:00"
